$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIO")

# Some Income/Cash-Flow figures were updated for the latest period, and a
# handful of most-recent-period (column J) figures became unavailable ("NA").

# Row 21 - Earnings Before Interest And Taxes
$ws.Range("J21").Value = "NA"

# Row 83 - Capital Expenditures
$ws.Range("J83").Value = "NA"

# Row 89 - Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 104100

# Row 91 - Dividends Paid
$ws.Range("D91").Value = -111300
$ws.Range("E91").Value = -141400
$ws.Range("F91").Value = -112000
$ws.Range("G91").Value = -121000
$ws.Range("H91").Value = -113000
$ws.Range("I91").Value = -152400
$ws.Range("J91").Value = -102900

# Row 94 - Total Cash Flows From Financing Activities
$ws.Range("J94").Value = "NA"

# Row 100 - Effect Of Exchange Rate Changes
$ws.Range("J100").Value = "NA"

# Row 101 - Change In Cash and Cash Equivalents
$ws.Range("J101").Value = "NA"

# Row 102 - (final row) Cash figure
$ws.Range("D102").Value = -72200

$wb.Save()
